$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 128.4548946666667
$ws.Range("H2").Value = 385.364684
$ws.Range("I2").Value = 0.2815548034715028
$ws.Range("J2").Value = 0.2815548034715028
$ws.Range("M2").Value = 23.641894
$ws.Range("N2").Value = 70.92568199999999
$ws.Range("O2").Value = 0.1609466983245457
$ws.Range("P2").Value = 0.1609466983245456
$ws.Range("Q2").Value = 3036.917003490498
$ws.Range("R2").Value = 27332.25303141449
$ws.Range("S2").Value = 0.0453153160161547
$ws.Range("T2").Value = 0.0453153160161547
$ws.Range("G3").Value = 128.4548946666667
$ws.Range("H3").Value = 385.364684
$ws.Range("I3").Value = 0.2815548034715028
$ws.Range("J3").Value = 0.2815548034715028
$ws.Range("O3").Value = 0.2271397161392734
$ws.Range("P3").Value = 0.2271397161392734
$ws.Range("Q3").Value = 4285.918712792654
$ws.Range("R3").Value = 38573.2684151339
$ws.Range("S3").Value = 0.06395227813816605
$ws.Range("T3").Value = 0.06395227813816605
$ws.Range("G4").Value = 128.4548946666667
$ws.Range("H4").Value = 385.364684
$ws.Range("I4").Value = 0.2815548034715028
$ws.Range("J4").Value = 0.2815548034715028
$ws.Range("M4").Value = 40.78693933333333
$ws.Range("N4").Value = 122.360818
$ws.Range("O4").Value = 0.2776648613881589
$ws.Range("P4").Value = 0.2776648613881589
$ws.Range("Q4").Value = 5239.281995839056
$ws.Range("R4").Value = 47153.53796255151
$ws.Range("S4").Value = 0.07817787547908515
$ws.Range("T4").Value = 0.07817787547908514
$ws.Range("G5").Value = 128.4548946666667
$ws.Range("H5").Value = 385.364684
$ws.Range("I5").Value = 0.2815548034715028
$ws.Range("J5").Value = 0.2815548034715028
$ws.Range("M5").Value = 10.35975466666667
$ws.Range("N5").Value = 31.079264
$ws.Range("O5").Value = 0.07052600392558668
$ws.Range("P5").Value = 0.07052600392558667
$ws.Range("Q5").Value = 1330.761194479175
$ws.Range("R5").Value = 11976.85075031258
$ws.Range("S5").Value = 0.01985693517489899
$ws.Range("T5").Value = 0.01985693517489899
$ws.Range("G6").Value = 128.4548946666667
$ws.Range("H6").Value = 385.364684
$ws.Range("I6").Value = 0.2815548034715028
$ws.Range("J6").Value = 0.2815548034715028
$ws.Range("M6").Value = 38.73894066666667
$ws.Range("N6").Value = 116.216822
$ws.Range("O6").Value = 0.2637227202224355
$ws.Range("P6").Value = 0.2637227202224354
$ws.Range("Q6").Value = 4976.206542834917
$ws.Range("R6").Value = 44785.85888551425
$ws.Range("S6").Value = 0.07425239866319794
$ws.Range("T6").Value = 0.07425239866319791
$ws.Range("H7").Value = 457.183265
$ws.Range("I7").Value = 0.3340268313936494
$ws.Range("J7").Value = 0.3340268313936494
$ws.Range("M7").Value = 23.641894
$ws.Range("N7").Value = 70.92568199999999
$ws.Range("O7").Value = 0.1609466983245457
$ws.Range("P7").Value = 0.1609466983245456
$ws.Range("Q7").Value = 3602.892763234636
$ws.Range("R7").Value = 32426.03486911173
$ws.Range("S7").Value = 0.05376051566461757
$ws.Range("T7").Value = 0.05376051566461756
$ws.Range("H8").Value = 457.183265
$ws.Range("I8").Value = 0.3340268313936494
$ws.Range("J8").Value = 0.3340268313936494
$ws.Range("O8").Value = 0.2271397161392734
$ws.Range("P8").Value = 0.2271397161392734
$ws.Range("Q8").Value = 5084.664973189767
$ws.Range("R8").Value = 45761.98475870791
$ws.Range("S8").Value = 0.07587075966565446
$ws.Range("T8").Value = 0.07587075966565446
$ws.Range("H9").Value = 457.183265
$ws.Range("I9").Value = 0.3340268313936494
$ws.Range("J9").Value = 0.3340268313936494
$ws.Range("M9").Value = 40.78693933333333
$ws.Range("N9").Value = 122.360818
$ws.Range("O9").Value = 0.2776648613881589
$ws.Range("P9").Value = 0.2776648613881589
$ws.Range("Q9").Value = 6215.702031256751
$ws.Range("R9").Value = 55941.31828131077
$ws.Range("S9").Value = 0.0927475138388436
$ws.Range("T9").Value = 0.09274751383884358
$ws.Range("H10").Value = 457.183265
$ws.Range("I10").Value = 0.3340268313936494
$ws.Range("J10").Value = 0.3340268313936494
$ws.Range("M10").Value = 10.35975466666667
$ws.Range("N10").Value = 31.079264
$ws.Range("O10").Value = 0.07052600392558668
$ws.Range("P10").Value = 0.07052600392558667
$ws.Range("Q10").Value = 1578.768821035218
$ws.Range("R10").Value = 14208.91938931696
$ws.Range("S10").Value = 0.0235575776221198
$ws.Range("T10").Value = 0.02355757762211979
$ws.Range("H11").Value = 457.183265
$ws.Range("I11").Value = 0.3340268313936494
$ws.Range("J11").Value = 0.3340268313936494
$ws.Range("M11").Value = 38.73894066666667
$ws.Range("N11").Value = 116.216822
$ws.Range("O11").Value = 0.2637227202224355
$ws.Range("P11").Value = 0.2637227202224354
$ws.Range("Q11").Value = 5903.598458875981
$ws.Range("R11").Value = 53132.38612988383
$ws.Range("S11").Value = 0.08809046460241403
$ws.Range("T11").Value = 0.088090464602414
$ws.Range("G12").Value = 70.798157
$ws.Range("H12").Value = 212.394471
$ws.Range("I12").Value = 0.1551794599342134
$ws.Range("J12").Value = 0.1551794599342134
$ws.Range("M12").Value = 23.641894
$ws.Range("N12").Value = 70.92568199999999
$ws.Range("O12").Value = 0.1609466983245457
$ws.Range("P12").Value = 0.1609466983245456
$ws.Range("Q12").Value = 1673.802523189358
$ws.Range("R12").Value = 15064.22270870422
$ws.Range("S12").Value = 0.02497562172419776
$ws.Range("T12").Value = 0.02497562172419776
$ws.Range("G13").Value = 70.798157
$ws.Range("H13").Value = 212.394471
$ws.Range("I13").Value = 0.1551794599342134
$ws.Range("J13").Value = 0.1551794599342134
$ws.Range("O13").Value = 0.2271397161392734
$ws.Range("P13").Value = 0.2271397161392734
$ws.Range("Q13").Value = 2362.192166401519
$ws.Range("R13").Value = 21259.72949761368
$ws.Range("S13").Value = 0.03524741848010298
$ws.Range("T13").Value = 0.03524741848010297
$ws.Range("G14").Value = 70.798157
$ws.Range("H14").Value = 212.394471
$ws.Range("I14").Value = 0.1551794599342134
$ws.Range("J14").Value = 0.1551794599342134
$ws.Range("M14").Value = 40.78693933333333
$ws.Range("N14").Value = 122.360818
$ws.Range("O14").Value = 0.2776648613881589
$ws.Range("P14").Value = 0.2776648613881589
$ws.Range("Q14").Value = 2887.640134470808
$ws.Range("R14").Value = 25988.76121023728
$ws.Range("S14").Value = 0.04308788323292272
$ws.Range("T14").Value = 0.04308788323292272
$ws.Range("G15").Value = 70.798157
$ws.Range("H15").Value = 212.394471
$ws.Range("I15").Value = 0.1551794599342134
$ws.Range("J15").Value = 0.1551794599342134
$ws.Range("M15").Value = 10.35975466666667
$ws.Range("N15").Value = 31.079264
$ws.Range("O15").Value = 0.07052600392558668
$ws.Range("P15").Value = 0.07052600392558667
$ws.Range("Q15").Value = 733.4515373721493
$ws.Range("R15").Value = 6601.063836349344
$ws.Range("S15").Value = 0.01094418720049076
$ws.Range("T15").Value = 0.01094418720049075
$ws.Range("G16").Value = 70.798157
$ws.Range("H16").Value = 212.394471
$ws.Range("I16").Value = 0.1551794599342134
$ws.Range("J16").Value = 0.1551794599342134
$ws.Range("M16").Value = 38.73894066666667
$ws.Range("N16").Value = 116.216822
$ws.Range("O16").Value = 0.2637227202224355
$ws.Range("P16").Value = 0.2637227202224354
$ws.Range("Q16").Value = 2742.645603332352
$ws.Range("R16").Value = 24683.81042999116
$ws.Range("S16").Value = 0.04092434929649919
$ws.Range("T16").Value = 0.04092434929649918
$ws.Range("G17").Value = 20.703408
$ws.Range("H17").Value = 62.110224
$ws.Range("I17").Value = 0.04537891674549766
$ws.Range("J17").Value = 0.04537891674549767
$ws.Range("M17").Value = 23.641894
$ws.Range("N17").Value = 70.92568199999999
$ws.Range("O17").Value = 0.1609466983245457
$ws.Range("P17").Value = 0.1609466983245456
$ws.Range("Q17").Value = 489.4677773747519
$ws.Range("R17").Value = 4405.209996372768
$ws.Range("S17").Value = 0.007303586823732285
$ws.Range("T17").Value = 0.007303586823732285
$ws.Range("G18").Value = 20.703408
$ws.Range("H18").Value = 62.110224
$ws.Range("I18").Value = 0.04537891674549766
$ws.Range("J18").Value = 0.04537891674549767
$ws.Range("O18").Value = 0.2271397161392734
$ws.Range("P18").Value = 0.2271397161392734
$ws.Range("Q18").Value = 690.772617081184
$ws.Range("R18").Value = 6216.953553730656
$ws.Range("S18").Value = 0.01030735426828006
$ws.Range("T18").Value = 0.01030735426828006
$ws.Range("G19").Value = 20.703408
$ws.Range("H19").Value = 62.110224
$ws.Range("I19").Value = 0.04537891674549766
$ws.Range("J19").Value = 0.04537891674549767
$ws.Range("M19").Value = 40.78693933333333
$ws.Range("N19").Value = 122.360818
$ws.Range("O19").Value = 0.2776648613881589
$ws.Range("P19").Value = 0.2776648613881589
$ws.Range("Q19").Value = 844.4286460892479
$ws.Range("R19").Value = 7599.857814803232
$ws.Range("S19").Value = 0.01260013062808341
$ws.Range("T19").Value = 0.01260013062808341
$ws.Range("G20").Value = 20.703408
$ws.Range("H20").Value = 62.110224
$ws.Range("I20").Value = 0.04537891674549766
$ws.Range("J20").Value = 0.04537891674549767
$ws.Range("M20").Value = 10.35975466666667
$ws.Range("N20").Value = 31.079264
$ws.Range("O20").Value = 0.07052600392558668
$ws.Range("P20").Value = 0.07052600392558667
$ws.Range("Q20").Value = 214.482227643904
$ws.Range("R20").Value = 1930.340048795136
$ws.Range("S20").Value = 0.003200393660531839
$ws.Range("T20").Value = 0.003200393660531839
$ws.Range("G21").Value = 20.703408
$ws.Range("H21").Value = 62.110224
$ws.Range("I21").Value = 0.04537891674549766
$ws.Range("J21").Value = 0.04537891674549767
$ws.Range("M21").Value = 38.73894066666667
$ws.Range("N21").Value = 116.216822
$ws.Range("O21").Value = 0.2637227202224355
$ws.Range("P21").Value = 0.2637227202224354
$ws.Range("Q21").Value = 802.0280941097921
$ws.Range("R21").Value = 7218.252846988129
$ws.Range("S21").Value = 0.01196745136487007
$ws.Range("T21").Value = 0.01196745136487007
$ws.Range("G22").Value = 83.88319133333333
$ws.Range("H22").Value = 251.649574
$ws.Range("I22").Value = 0.1838599884551367
$ws.Range("J22").Value = 0.1838599884551367
$ws.Range("M22").Value = 23.641894
$ws.Range("N22").Value = 70.92568199999999
$ws.Range("O22").Value = 0.1609466983245457
$ws.Range("P22").Value = 0.1609466983245456
$ws.Range("Q22").Value = 1983.157517884385
$ws.Range("R22").Value = 17848.41766095947
$ws.Range("S22").Value = 0.02959165809584333
$ws.Range("T22").Value = 0.02959165809584333
$ws.Range("G23").Value = 83.88319133333333
$ws.Range("H23").Value = 251.649574
$ws.Range("I23").Value = 0.1838599884551367
$ws.Range("J23").Value = 0.1838599884551367
$ws.Range("O23").Value = 0.2271397161392734
$ws.Range("P23").Value = 0.2271397161392734
$ws.Range("Q23").Value = 2798.776491602172
$ws.Range("R23").Value = 25188.98842441956
$ws.Range("S23").Value = 0.04176190558706983
$ws.Range("T23").Value = 0.04176190558706983
$ws.Range("G24").Value = 83.88319133333333
$ws.Range("H24").Value = 251.649574
$ws.Range("I24").Value = 0.1838599884551367
$ws.Range("J24").Value = 0.1838599884551367
$ws.Range("M24").Value = 40.78693933333333
$ws.Range("N24").Value = 122.360818
$ws.Range("O24").Value = 0.2776648613881589
$ws.Range("P24").Value = 0.2776648613881589
$ws.Range("Q24").Value = 3421.338635999059
$ws.Range("R24").Value = 30792.04772399153
$ws.Range("S24").Value = 0.05105145820922403
$ws.Range("T24").Value = 0.05105145820922402
$ws.Range("G25").Value = 83.88319133333333
$ws.Range("H25").Value = 251.649574
$ws.Range("I25").Value = 0.1838599884551367
$ws.Range("J25").Value = 0.1838599884551367
$ws.Range("M25").Value = 10.35975466666667
$ws.Range("N25").Value = 31.079264
$ws.Range("O25").Value = 0.07052600392558668
$ws.Range("P25").Value = 0.07052600392558667
$ws.Range("Q25").Value = 869.0092828703928
$ws.Range("R25").Value = 7821.083545833536
$ws.Range("S25").Value = 0.01296691026754529
$ws.Range("T25").Value = 0.01296691026754529
$ws.Range("G26").Value = 83.88319133333333
$ws.Range("H26").Value = 251.649574
$ws.Range("I26").Value = 0.1838599884551367
$ws.Range("J26").Value = 0.1838599884551367
$ws.Range("M26").Value = 38.73894066666667
$ws.Range("N26").Value = 116.216822
$ws.Range("O26").Value = 0.2637227202224355
$ws.Range("P26").Value = 0.2637227202224354
$ws.Range("Q26").Value = 3249.545971992648
$ws.Range("R26").Value = 29245.91374793383
$ws.Range("S26").Value = 0.04848805629545423
$ws.Range("T26").Value = 0.0484880562954542
